# Natalia Bernal, creación de dasboard
#
# Updates the "BODEGA" (warehouse) column so each count line references its
# own named warehouse ("Bodega 1".."Bodega 5") instead of the generic
# "BODEGA " / "BODEGA SEGUNDARIA" placeholders, fixes the OPERARIO count in
# K3, and leaves the cursor selected around C17:C18 like the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Bodega 1"
$ws.Range("A3").Value = "Bodega 2"
$ws.Range("A4").Value = "Bodega 3"
$ws.Range("A5").Value = "Bodega 4"
$ws.Range("A6").Value = "Bodega 5"

$ws.Range("K3").Value = 1

$ws.Range("C17:C18").Select()
